# Add new "flash access" entries to Sheet1 (zigbee_xbr flash分配 table)
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Insert 3 new rows right before the current row 9 ("resetbtcnt"), pushing
# the two trailing rows (resetbtcnt, zigbee_join_cnt) down to rows 12-13.
$insertRange = $ws.Range("A9:E11")
$insertRange.EntireRow.Insert()

# Row 9: Linkage_flag
$ws.Cells.Item(9, 1).Value = 8
$ws.Cells.Item(9, 2).Value = "Linkage_flag"
$ws.Cells.Item(9, 3).Value = "u8"
$ws.Cells.Item(9, 4).Value = "联动标志位"
$ws.Cells.Item(9, 5).Value = "0x2F08"

# Row 10: SWITCHflag2
$ws.Cells.Item(10, 1).Value = 9
$ws.Cells.Item(10, 2).Value = "SWITCHflag2"
$ws.Cells.Item(10, 3).Value = "u8"
$ws.Cells.Item(10, 4).Value = "开关灯"
$ws.Cells.Item(10, 5).Value = "0x2F09"

# Row 11: all_day_micro_light_enable
$ws.Cells.Item(11, 1).Value = 10
$ws.Cells.Item(11, 2).Value = "all_day_micro_light_enable"
$ws.Cells.Item(11, 3).Value = "u8"
$ws.Cells.Item(11, 4).Value = "全天伴亮"
$ws.Cells.Item(11, 5).Value = "0x2F0A"

# Update sequence numbers on the rows that got pushed down
$ws.Cells.Item(12, 1).Value = 11
$ws.Cells.Item(13, 1).Value = 12

# Apply same style/border as the rest of the data rows to the new cells
$styleSource = $ws.Range("A8:E8")
$styleSource.Copy()
$ws.Range("A9:E11").PasteSpecial(-4122) # xlPasteFormats
$excel.CutCopyMode = $false

# Widen column B to fit the new longer identifier names
$ws.Columns.Item(2).ColumnWidth = 27.142857142857142

# Update selection to match final author state
$ws.Range("B18").Select()
